$d = $word.ActiveDocument

# The last body paragraph contains the "entrepreneur / land request" text,
# split across several runs. Collapse that entire paragraph's text down to
# a single space, while keeping the paragraph's own formatting/rPr.
$lastTextParaIndex = $d.Paragraphs.Count - 1
$p = $d.Paragraphs.Item($lastTextParaIndex)
$r = $p.Range
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Text = " "

# Remove the now-stray trailing empty paragraph that followed it, by
# deleting from the end of the text paragraph's own paragraph mark through
# the end of the following (empty) paragraph's mark. This merges the two
# paragraph marks into one, leaving the text paragraph as the last one.
$p = $d.Paragraphs.Item($lastTextParaIndex)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRange = $d.Range($p.Range.End - 1, $lastPara.Range.End)
$delRange.Delete()
